$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status went from "Ready for handoff" to "In Translation" everywhere it appears
# (the Overview rollup columns for each locale, plus each locale sheet's own
# Status column).
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# Narrow the Status-related columns (the stored col width is ColumnWidth + 5/6,
# and this engine snaps ColumnWidth to the nearest 1/6 character, so 12.5 is
# the closest achievable setting to the authored 13.4101845877511 width).
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
